# Fruta / hortaliza, semanal
# Insert a new weekly record at row 338 (pushing the existing rows 338-448
# down to 339-449), then populate the new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 338:448 down to 339:449 by inserting a new row at 338.
$ws.Rows("338:338").Insert()

# Populate the newly inserted row 338 with the new weekly record.
$ws.Range("A338").Value = 5
$ws.Range("B338").Value = "Macroferia Regional de Talca"
$ws.Range("C338").Value = "Maule"
$ws.Range("D338").Value = 44627
$ws.Range("E338").Value = 7
$ws.Range("F338").Value = "Fruta"
$ws.Range("G338").Value = 100102
$ws.Range("H338").Value = "Cítricos"
$ws.Range("I338").Value = 100102005
$ws.Range("J338").Value = "Naranja"
$ws.Range("K338").Value = "Valencia"
$ws.Range("L338").Value = "Primera"
$ws.Range("M338").Value = 390
$ws.Range("N338").Value = 9000
$ws.Range("O338").Value = 10000
$ws.Range("P338").Value = 9513
$ws.Range("Q338").Value = "$/bandeja 15 kilos granel"
$ws.Range("R338").Value = "Región de O'Higgins"
$ws.Range("S338").Value = 634
$ws.Range("T338").Value = 15
